# Euro 2016 Game_Score id review
# The id sequence for the "2016" sheet's game_score insert statements
# started from '2008'!A192+1. It should instead continue from the
# '2012' sheet's last id ('2012'!A188+1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2016")

$ws.Range("A81").Formula = "='2012'!A188+1"

$wb.Application.Calculate()
